$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing existing rows 6-13 down to 7-14
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with its values
$ws.Range("A6").Value = "gre"
$ws.Range("B6").Value = "g"

# Update the selection to match the saved workbook state
$ws.Range("F8").Select()
